$d = $word.ActiveDocument

# --- Edit 1: "Steel 2023" -> "Steel 2024" ---------------------------------
# Locate the exact occurrence of "Steel 2023" in the document body.
$find1 = $d.Content
$find1.Find.Execute("Steel 2023", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $find1.Start
$end1 = $find1.End

# The final digit ("3") is the last character of the match.
$digitStart1 = $end1 - 1
$digitEnd1 = $end1

# Step 1: split the trailing digit off into its own run without disturbing
# neighbouring runs -- a pure formatting round-trip (no text change) does not
# trigger the engine's "merge adjacent same-format runs" normalisation.
$r1a = $d.Range($digitStart1, $digitEnd1)
$r1a.Font.Name = "Arial"
$r1b = $d.Range($digitStart1, $digitEnd1)
$r1b.Font.Name = "Verdana"

# Step 2: change the isolated run's text while its formatting still differs
# from its neighbours, so the content edit cannot re-merge it with them.
$r1c = $d.Range($digitStart1, $digitEnd1)
$r1c.Font.Name = "Arial"
$r1c.Text = "4"

# Step 3: restore the run's formatting to match its neighbours again. This is
# a pure formatting write, so it will not re-trigger a run merge.
$r1d = $d.Range($digitStart1, $digitEnd1)
$r1d.Font.Name = "Verdana"

# --- Edit 2: "2022" -> "2023" ----------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("2022", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $find2.Start
$end2 = $find2.End

$digitStart2 = $end2 - 1
$digitEnd2 = $end2

$r2a = $d.Range($digitStart2, $digitEnd2)
$r2a.Font.Name = "Arial"
$r2b = $d.Range($digitStart2, $digitEnd2)
$r2b.Font.Name = "Verdana"

$r2c = $d.Range($digitStart2, $digitEnd2)
$r2c.Font.Name = "Arial"
$r2c.Text = "3"

$r2d = $d.Range($digitStart2, $digitEnd2)
$r2d.Font.Name = "Verdana"
